$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("M2").ClearContents()
$ws.Range("N2").ClearContents()
$ws.Range("H9").Value = 537.5
$ws.Range("I9").Value = 537.5
$ws.Range("K9").Value = 537.5
$ws.Range("M9").Value = -368.5
$ws.Range("H88").Value = 7757.3335
$ws.Range("I88").Value = 7499.5
$ws.Range("J88").Value = 8273
$ws.Range("K88").Value = 7499.5
$ws.Range("L88").Value = 8273
$ws.Range("M88").Value = -7093.5
$ws.Range("N88").Value = -9085
$ws.Range("H91").Value = 7757.3335
$ws.Range("I91").Value = 7499.5
$ws.Range("J91").Value = 8273
$ws.Range("K91").Value = 7499.5
$ws.Range("L91").Value = 8273
$ws.Range("M91").Value = -6095.5
$ws.Range("N91").Value = -11081
$ws.Range("H98").Value = 1000
$ws.Range("I98").Value = 1000
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 1000
$ws.Range("L98").Value = 0
$ws.Range("M98").Value = 498
$ws.Range("N98").ClearContents()
$ws.Range("H100").Value = 4900
$ws.Range("I100").Value = 4900
$ws.Range("K100").Value = 4900
$ws.Range("M100").Value = -4359
$ws.Range("H113").Value = 11816.667
$ws.Range("I113").Value = 5000
$ws.Range("J113").Value = 18633.334
$ws.Range("K113").Value = 5000
$ws.Range("L113").Value = 18633.334
$ws.Range("M113").Value = -1746
$ws.Range("N113").Value = -25141.334
$ws.Range("H122").Value = 1000
$ws.Range("I122").Value = 1000
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 3000
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -550
$ws.Range("N122").ClearContents()
$ws.Range("H132").Value = 4612.6665
$ws.Range("J132").Value = 7063.25
$ws.Range("L132").Value = 21189.75
$ws.Range("N132").Value = -26249.75
$ws.Range("H137").Value = 3068.4211
$ws.Range("I137").Value = 2575
$ws.Range("K137").Value = 7725
$ws.Range("M137").Value = -5175

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 5000
$ws.Range("I61").Value = 5000
$ws.Range("K61").Value = 5000
$ws.Range("M61").Value = -4788
$ws.Range("H74").Value = 3125.6
$ws.Range("I74").Value = 3125.6
$ws.Range("K74").Value = 3125.6
$ws.Range("M74").Value = -2251.6
$ws.Range("H77").Value = 3125.6
$ws.Range("I77").Value = 3125.6
$ws.Range("K77").Value = 15628
$ws.Range("M77").Value = -11260
$ws.Range("H136").Value = 5000
$ws.Range("I136").Value = 5000
$ws.Range("K136").Value = 15000
$ws.Range("M136").Value = -12450

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 634.36365
$ws.Range("I80").Value = 547.8570999999999
$ws.Range("K80").Value = 547.8570999999999
$ws.Range("M80").Value = 450.1429000000001
$ws.Range("H83").Value = 634.36365
$ws.Range("I83").Value = 547.8570999999999
$ws.Range("K83").Value = 2739.2855
$ws.Range("M83").Value = 2252.7145
$ws.Range("H86").Value = 3510.6
$ws.Range("I86").Value = 1945
$ws.Range("J86").Value = 5299.857
$ws.Range("K86").Value = 1945
$ws.Range("L86").Value = 5299.857
$ws.Range("M86").Value = -822
$ws.Range("N86").Value = -7545.857
$ws.Range("H89").Value = 3510.6
$ws.Range("I89").Value = 1945
$ws.Range("J89").Value = 5299.857
$ws.Range("K89").Value = 9725
$ws.Range("L89").Value = 26499.285
$ws.Range("M89").Value = -4109
$ws.Range("N89").Value = -37731.285

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2000
$ws.Range("I16").Value = 2000
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 2000
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -1713
$ws.Range("N16").ClearContents()
$ws.Range("H22").Value = 767.3333
$ws.Range("I22").Value = 650.5
$ws.Range("K22").Value = 650.5
$ws.Range("M22").Value = -300.5
$ws.Range("H50").Value = 28571.428
$ws.Range("I50").Value = 26666.666
$ws.Range("K50").Value = 26666.666
$ws.Range("M50").Value = -26041.666
$ws.Range("H58").Value = 2487
$ws.Range("I58").Value = 2487
$ws.Range("K58").Value = 2487
$ws.Range("M58").Value = -2284
$ws.Range("H99").Value = 0
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("L99").Value = 0
$ws.Range("M99").ClearContents()
$ws.Range("N99").ClearContents()
$ws.Range("H113").Value = 2000
$ws.Range("I113").Value = 2000
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 2000
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 170
$ws.Range("N113").ClearContents()
$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("M126").ClearContents()
$ws.Range("N126").ClearContents()
$ws.Range("H136").Value = 2487
$ws.Range("I136").Value = 2487
$ws.Range("K136").Value = 7461
$ws.Range("M136").Value = -4911
$ws.Range("H141").Value = 117446
$ws.Range("J141").Value = 117446
$ws.Range("L141").Value = 117446
$ws.Range("N141").Value = -127806

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 527.8
$ws.Range("J113").Value = 549.5
$ws.Range("L113").Value = 1648.5
$ws.Range("N113").Value = -5988.5
$ws.Range("H121").Value = 1461.2
$ws.Range("I121").Value = 665
$ws.Range("J121").Value = 1660.25
$ws.Range("K121").Value = 1995
$ws.Range("L121").Value = 4980.75
$ws.Range("M121").Value = -685
$ws.Range("N121").Value = -7600.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 8669.5
$ws.Range("I43").Value = 2339
$ws.Range("K43").Value = 2339
$ws.Range("M43").Value = -2188
$ws.Range("H46").Value = 17500
$ws.Range("I46").Value = 10000
$ws.Range("K46").Value = 10000
$ws.Range("M46").Value = -9844
$ws.Range("H57").Value = 12138.75
$ws.Range("I57").Value = 7851.6665
$ws.Range("K57").Value = 7851.6665
$ws.Range("M57").Value = -7031.6665
$ws.Range("H132").Value = 6388.5
$ws.Range("I132").Value = 6388.5
$ws.Range("K132").Value = 19165.5
$ws.Range("M132").Value = -16635.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1228.125
$ws.Range("I22").Value = 473.33334
$ws.Range("J22").Value = 1681
$ws.Range("K22").Value = 473.33334
$ws.Range("L22").Value = 1681
$ws.Range("M22").Value = -178.33334
$ws.Range("N22").Value = -2271
$ws.Range("H27").Value = 1228.125
$ws.Range("I27").Value = 473.33334
$ws.Range("J27").Value = 1681
$ws.Range("K27").Value = 473.33334
$ws.Range("L27").Value = 1681
$ws.Range("M27").Value = -366.33334
$ws.Range("N27").Value = -1895
$ws.Range("H40").Value = 7500
$ws.Range("I40").Value = 7500
$ws.Range("J40").Value = 7500
$ws.Range("K40").Value = 7500
$ws.Range("L40").Value = 7500
$ws.Range("M40").Value = -7364
$ws.Range("N40").Value = -7772
$ws.Range("H132").Value = 32249.666
$ws.Range("I132").Value = 32249.666
$ws.Range("K132").Value = 96748.99800000001
$ws.Range("M132").Value = -94218.99800000001
$ws.Range("H136").Value = 2467.6
$ws.Range("I136").Value = 2459.5
$ws.Range("J136").Value = 2500
$ws.Range("K136").Value = 7378.5
$ws.Range("L136").Value = 7500
$ws.Range("M136").Value = -4828.5
$ws.Range("N136").Value = -12600

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 546.9167
$ws.Range("I107").Value = 418.22223
$ws.Range("K107").Value = 1254.66669
$ws.Range("M107").Value = 665.33331
$ws.Range("H136").Value = 2609.625
$ws.Range("I136").Value = 2650.2666
$ws.Range("J136").Value = 2000
$ws.Range("K136").Value = 7950.7998
$ws.Range("L136").Value = 6000
$ws.Range("M136").Value = -5400.7998
$ws.Range("N136").Value = -11100
